$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.052.13'
$ws.Range('E2').Value = '  -1.30%  '
$ws.Range('D3').Value = '2.299.65'
$ws.Range('E3').Value = '  -2.06%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'313.47"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.96%  '
$ws.Range('D6').Value = "'106.60"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.02%  '
$ws.Range('D7').Value = "'0.629"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.87%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -1.39%  '
$ws.Range('D10').Value = "'40.35"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.00%  '
$ws.Range('D11').Value = "'0.0914"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('D12').Value = "'8.32"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.96%  '
$ws.Range('E13').Value = '  +0.08%  '
$ws.Range('D14').Value = "'0.975"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.96%  '
$ws.Range('D15').Value = "'15.58"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.11%  '
$ws.Range('D16').Value = '2.644.52'
$ws.Range('E16').Value = '  -2.27%  '
$ws.Range('D17').Value = '2.315.90'
$ws.Range('E17').Value = '  -1.57%  '
$ws.Range('D18').Value = '42.095.16'
$ws.Range('E18').Value = '  -1.54%  '
$ws.Range('D19').Value = "'7.55"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.35%  '
$ws.Range('E20').Value = '  -1.77%  '
$ws.Range('D21').Value = "'73.21"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.61%  '
$ws.Range('E22').Value = '  -4.08%  '
$ws.Range('D23').Value = "'257.40"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.34%  '
$ws.Range('D24').Value = "'2.32"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('E25').Value = '  -4.26%  '
$ws.Range('E26').Value = '  +0.39%  '
$ws.Range('E27').Value = '  -3.86%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = "'22.88"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.35%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = "'2.23"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.01%  '
$ws.Range('D30').Value = "'166.36"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.68%  '
$ws.Range('D31').Value = "'35.78"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.52%  '
$ws.Range('E32').Value = '  -0.34%  '
$ws.Range('E33').Value = '  -7.12%  '
$ws.Range('E34').Value = '  -7.17%  '
$ws.Range('E35').Value = '  +6.39%  '
$ws.Range('E36').Value = '  -2.15%  '
$ws.Range('E37').Value = '  +0.92%  '
$ws.Range('E38').Value = '  -1.41%  '
$ws.Range('D39').Value = "'2.92"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.83%  '
$ws.Range('E40').Value = '  -4.11%  '
$ws.Range('E41').Value = '  +2.42%  '
$ws.Range('D42').Value = "'72.28"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.66%  '
$ws.Range('D43').Value = "'98.07"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.27%  '
$ws.Range('E44').Value = '  -3.00%  '
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('D46').Value = "'12.37"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.28%  '
$ws.Range('D47').Value = "'113.16"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.97%  '
$ws.Range('D48').Value = "'9.11"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.32%  '
$ws.Range('B49').Value = 'ordi'
$ws.Range('C49').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D49').Value = "'76.44"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.88%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = "'5.33"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.26%  '
$ws.Range('D51').Value = "'1.27"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.11%  '
